$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 11983.343
$ws.Range("I112").Value = 550
$ws.Range("J112").Value = 12305.408
$ws.Range("K112").Value = 1650
$ws.Range("L112").Value = 36916.224
$ws.Range("M112").Value = -542
$ws.Range("N112").Value = -39132.224
$ws.Range("H137").Value = 1463.7297
$ws.Range("I137").Value = 1255.1613
$ws.Range("J137").Value = 2541.3333
$ws.Range("K137").Value = 3765.4839
$ws.Range("L137").Value = 7623.999899999999
$ws.Range("M137").Value = -1215.4839
$ws.Range("N137").Value = -12723.9999
$ws.Range("H141").Value = 3216.644
$ws.Range("I141").Value = 1309.5385
$ws.Range("J141").Value = 6935.5
$ws.Range("K141").Value = 3928.6155
$ws.Range("L141").Value = 20806.5
$ws.Range("M141").Value = 1251.3845
$ws.Range("N141").Value = -31166.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1272.5217
$ws.Range("I61").Value = 969.1111
$ws.Range("J61").Value = 2364.8
$ws.Range("K61").Value = 969.1111
$ws.Range("L61").Value = 2364.8
$ws.Range("M61").Value = -757.1111
$ws.Range("N61").Value = -2788.8
$ws.Range("H74").Value = 1260.5918
$ws.Range("I74").Value = 1262.9773
$ws.Range("K74").Value = 1262.9773
$ws.Range("M74").Value = -388.9773
$ws.Range("H77").Value = 1260.5918
$ws.Range("I77").Value = 1262.9773
$ws.Range("K77").Value = 6314.886500000001
$ws.Range("M77").Value = -1946.886500000001
$ws.Range("H136").Value = 1272.5217
$ws.Range("I136").Value = 969.1111
$ws.Range("J136").Value = 2364.8
$ws.Range("K136").Value = 2907.3333
$ws.Range("L136").Value = 7094.400000000001
$ws.Range("M136").Value = -357.3332999999998
$ws.Range("N136").Value = -12194.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3090419.8
$ws.Range("I134").Value = 1255.8572
$ws.Range("K134").Value = 3767.5716
$ws.Range("M134").Value = -1232.5716

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1268.4138
$ws.Range("I31").Value = 1000.8095
$ws.Range("J31").Value = 1518.1777
$ws.Range("K31").Value = 1000.8095
$ws.Range("L31").Value = 1518.1777
$ws.Range("M31").Value = -705.8095
$ws.Range("N31").Value = -2108.1777
$ws.Range("H34").Value = 1268.4138
$ws.Range("I34").Value = 1000.8095
$ws.Range("J34").Value = 1518.1777
$ws.Range("K34").Value = 1000.8095
$ws.Range("L34").Value = 1518.1777
$ws.Range("M34").Value = -798.8095
$ws.Range("N34").Value = -1922.1777
$ws.Range("H58").Value = 15152310
$ws.Range("I58").Value = 30303742
$ws.Range("J58").Value = 879.30304
$ws.Range("K58").Value = 30303742
$ws.Range("L58").Value = 879.30304
$ws.Range("M58").Value = -30303539
$ws.Range("N58").Value = -1285.30304
$ws.Range("H132").Value = 6803825
$ws.Range("I132").Value = 951.3714
$ws.Range("J132").Value = 23811010
$ws.Range("K132").Value = 2854.1142
$ws.Range("L132").Value = 71433030
$ws.Range("M132").Value = -324.1142
$ws.Range("N132").Value = -71438090
$ws.Range("H134").Value = 969.7586
$ws.Range("I134").Value = 881
$ws.Range("J134").Value = 1167
$ws.Range("K134").Value = 2643
$ws.Range("L134").Value = 3501
$ws.Range("M134").Value = -108
$ws.Range("N134").Value = -8571
$ws.Range("H136").Value = 15152310
$ws.Range("I136").Value = 30303742
$ws.Range("J136").Value = 879.30304
$ws.Range("K136").Value = 90911226
$ws.Range("L136").Value = 2637.90912
$ws.Range("M136").Value = -90908676
$ws.Range("N136").Value = -7737.90912

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 13334395
$ws.Range("I129").Value = 841.46155
$ws.Range("J129").Value = 27779078
$ws.Range("K129").Value = 2524.38465
$ws.Range("L129").Value = 83337234
$ws.Range("M129").Value = 2475.61535
$ws.Range("N129").Value = -83347234
$ws.Range("H130").Value = 50001376
$ws.Range("I130").Value = 125000104
$ws.Range("J130").Value = 2222.1667
$ws.Range("K130").Value = 375000312
$ws.Range("L130").Value = 6666.500100000001
$ws.Range("M130").Value = -374995292
$ws.Range("N130").Value = -16706.5001
$ws.Range("H131").Value = 13568357
$ws.Range("I131").Value = 51283230
$ws.Range("J131").Value = 6852009
$ws.Range("K131").Value = 153849690
$ws.Range("L131").Value = 20556027
$ws.Range("M131").Value = -153844650
$ws.Range("N131").Value = -20566107
$ws.Range("H133").Value = 47622480
$ws.Range("J133").Value = 5000
$ws.Range("L133").Value = 15000
$ws.Range("N133").Value = -25120
$ws.Range("H134").Value = 38464500
$ws.Range("I134").Value = 83334380
$ws.Range("J134").Value = 4598.2856
$ws.Range("K134").Value = 250003140
$ws.Range("L134").Value = 13794.8568
$ws.Range("M134").Value = -249998070
$ws.Range("N134").Value = -23934.8568
$ws.Range("H136").Value = 46299772
$ws.Range("I136").Value = 89287050
$ws.Range("J136").Value = 5783.6924
$ws.Range("K136").Value = 267861150
$ws.Range("L136").Value = 17351.0772
$ws.Range("M136").Value = -267856050
$ws.Range("N136").Value = -27551.0772
$ws.Range("H137").Value = 32052948
$ws.Range("I137").Value = 35715324
$ws.Range("J137").Value = 27780176
$ws.Range("K137").Value = 107145972
$ws.Range("L137").Value = 83340528
$ws.Range("M137").Value = -107140872
$ws.Range("N137").Value = -83350728
$ws.Range("H138").Value = 23189910
$ws.Range("I138").Value = 31373624
$ws.Range("K138").Value = 94120872
$ws.Range("M138").Value = -94115732
$ws.Range("H139").Value = 13439219
$ws.Range("I139").Value = 41667720
$ws.Range("J139").Value = 410679.47
$ws.Range("K139").Value = 125003160
$ws.Range("L139").Value = 1232038.41
$ws.Range("M139").Value = -124998020
$ws.Range("N139").Value = -1242318.41
$ws.Range("H140").Value = 23686316
$ws.Range("I140").Value = 34617080
$ws.Range("J140").Value = 2999.8333
$ws.Range("K140").Value = 103851240
$ws.Range("L140").Value = 8999.499899999999
$ws.Range("M140").Value = -103846060
$ws.Range("N140").Value = -19359.4999
$ws.Range("H141").Value = 40002492
$ws.Range("I141").Value = 47621016
$ws.Range("J141").Value = 5249.5
$ws.Range("K141").Value = 142863048
$ws.Range("L141").Value = 15748.5
$ws.Range("M141").Value = -142857868
$ws.Range("N141").Value = -26108.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 19613704
$ws.Range("I132").Value = 29413386
$ws.Range("J132").Value = 14335.471
$ws.Range("K132").Value = 88240158
$ws.Range("L132").Value = 43006.413
$ws.Range("M132").Value = -88237628
$ws.Range("N132").Value = -48066.413
$ws.Range("H136").Value = 30076910
$ws.Range("I136").Value = 5293021
$ws.Range("J136").Value = 90910090
$ws.Range("K136").Value = 15879063
$ws.Range("L136").Value = 272730270
$ws.Range("M136").Value = -15876513
$ws.Range("N136").Value = -272735370

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6426278
$ws.Range("I132").Value = 19193.21
$ws.Range("J132").Value = 23816938
$ws.Range("K132").Value = 57579.63
$ws.Range("L132").Value = 71450814
$ws.Range("M132").Value = -55049.63
$ws.Range("N132").Value = -71455874
$ws.Range("H136").Value = 10643363
$ws.Range("I136").Value = 12825776
$ws.Range("J136").Value = 4099.375
$ws.Range("K136").Value = 38477328
$ws.Range("L136").Value = 12298.125
$ws.Range("M136").Value = -38474778
$ws.Range("N136").Value = -17398.125
